$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "New Blank Display" / "ND" shortcut row.
$ws.Rows.Item(88).Delete()
